# ST-246: split the title run into "THÊM " + "CHỨC VỤ", and mark the
# two picture runs as NoProofing (<w:noProof/>).

$d = $word.ActiveDocument

# --- 1. Split "THÊM CHƯƠNG TRÌNH KHUYẾN MÃI" into two runs: ---
#        "THÊM " (kept) + "CHỨC VỤ" (new text), same run formatting.
$titleRange = $d.Range(0, 28)
if ($titleRange.Text -eq "THÊM CHƯƠNG TRÌNH KHUYẾN MÃI") {
    $xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">THÊM </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>CHỨC VỤ</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $titleRange.InsertXML($xml)
}

# --- 2. Mark the picture runs (the two <w:drawing> inline pictures) as
#        NoProofing, i.e. add <w:noProof/> to their run properties. The
#        three text-box drawings already carry <w:noProof/> and are left
#        untouched.
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = $true
}
